$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.35038893185310371
$ws.Range("A2").Value = -0.042494950880548998
$ws.Range("A3").Value = -0.0089999998425920325
$ws.Range("A4").Value = 0.061998139567201349
$ws.Range("A5").Value = -0.0059999998453923453
$ws.Range("A6").Value = -0.0059999998395205978
$ws.Range("A7").Value = -0.019999999805820678
$ws.Range("A8").Value = 0.019459576083664665
$ws.Range("A9").Value = -0.0059999998345290351
$ws.Range("A10").Value = -0.0059999998316868641
$ws.Range("A11").Value = -0.004499999834866486
$ws.Range("A12").Value = -0.0059999998305428903
$ws.Range("A13").Value = -0.005999999827574598
$ws.Range("A14").Value = -0.01199999981242339
$ws.Range("A15").Value = -0.0059999998259359089
$ws.Range("A16").Value = -0.0059999998252684428
$ws.Range("A17").Value = -0.0059999998243993602
$ws.Range("A18").Value = -0.0089999998170604556
$ws.Range("A19").Value = -0.0089999998463388131
$ws.Range("A20").Value = -0.0089999998422474192
$ws.Range("A21").Value = -0.05425762947912105
$ws.Range("A22").Value = -0.0089999998411296467
$ws.Range("A23").Value = -0.0089999998410439375
$ws.Range("A24").Value = -0.041999999760084172
$ws.Range("A25").Value = -0.041999999758715489
$ws.Range("A26").Value = -0.0059999998387958442
$ws.Range("A27").Value = -0.0059999998378046371
$ws.Range("A28").Value = -0.0059999998340067862
$ws.Range("A29").Value = -0.01199999981711386
$ws.Range("A30").Value = -0.01999999979681899
$ws.Range("A31").Value = -0.014999999806802222
$ws.Range("A32").Value = -0.020999999792231883
$ws.Range("A33").Value = -0.0059999998276758504
